$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3758.5667
$ws.Range("I98").Value = 2633.353
$ws.Range("J98").Value = 5230
$ws.Range("K98").Value = 2633.353
$ws.Range("L98").Value = 5230
$ws.Range("M98").Value = -1135.353
$ws.Range("N98").Value = -8226
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H116").Value = 4496.222
$ws.Range("I116").Value = 4994.5
$ws.Range("J116").Value = 3499.6667
$ws.Range("K116").Value = 4994.5
$ws.Range("L116").Value = 3499.6667
$ws.Range("M116").Value = -1552.5
$ws.Range("N116").Value = -10383.6667
$ws.Range("H122").Value = 3758.5667
$ws.Range("I122").Value = 2633.353
$ws.Range("J122").Value = 5230
$ws.Range("K122").Value = 7900.059
$ws.Range("L122").Value = 15690
$ws.Range("M122").Value = -5450.059
$ws.Range("N122").Value = -20590
$ws.Range("H132").Value = 1661.6316
$ws.Range("I132").Value = 1411.375
$ws.Range("K132").Value = 4234.125
$ws.Range("M132").Value = -1704.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 751.5357
$ws.Range("I2").Value = 781.61536
$ws.Range("K2").Value = 781.61536
$ws.Range("M2").Value = -668.61536
$ws.Range("I8").Value = 7400
$ws.Range("J8").Value = 14166
$ws.Range("K8").Value = 7400
$ws.Range("L8").Value = 14166
$ws.Range("M8").Value = -7256
$ws.Range("N8").Value = -14454
$ws.Range("H32").Value = 219959.12
$ws.Range("I32").Value = 257461.53
$ws.Range("K32").Value = 257461.53
$ws.Range("M32").Value = -257174.53
$ws.Range("H45").Value = 40714.69
$ws.Range("I45").Value = 64813.125
$ws.Range("K45").Value = 64813.125
$ws.Range("M45").Value = -64436.125
$ws.Range("H74").Value = 466729.25
$ws.Range("I74").Value = 1420.4375
$ws.Range("K74").Value = 1420.4375
$ws.Range("M74").Value = -546.4375
$ws.Range("H77").Value = 466729.25
$ws.Range("I77").Value = 1420.4375
$ws.Range("K77").Value = 7102.1875
$ws.Range("M77").Value = -2734.1875
$ws.Range("H88").Value = 1866.8
$ws.Range("J88").Value = 1811.75
$ws.Range("L88").Value = 1811.75
$ws.Range("N88").Value = -2623.75
$ws.Range("H91").Value = 1866.8
$ws.Range("J91").Value = 1811.75
$ws.Range("L91").Value = 1811.75
$ws.Range("N91").Value = -4619.75
$ws.Range("H116").Value = 751.5357
$ws.Range("I116").Value = 781.61536
$ws.Range("K116").Value = 781.61536
$ws.Range("M116").Value = 1512.38464
$ws.Range("H119").Value = 58939
$ws.Range("J119").Value = 58939
$ws.Range("L119").Value = 58939
$ws.Range("N119").Value = -68615
$ws.Range("H122").Value = 2755.75
$ws.Range("I122").Value = 2435.1428
$ws.Range("K122").Value = 7305.428400000001
$ws.Range("M122").Value = -4855.428400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 751.5357
$ws.Range("I3").Value = 781.61536
$ws.Range("K3").Value = 781.61536
$ws.Range("M3").Value = -667.61536
$ws.Range("H20").Value = 804.931
$ws.Range("I20").Value = 769.45
$ws.Range("J20").Value = 883.7778
$ws.Range("K20").Value = 769.45
$ws.Range("L20").Value = 883.7778
$ws.Range("M20").Value = -522.45
$ws.Range("N20").Value = -1377.7778
$ws.Range("H50").Value = 64999.5
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H86").Value = 2466.5625
$ws.Range("I86").Value = 1547.25
$ws.Range("J86").Value = 5224.5
$ws.Range("K86").Value = 1547.25
$ws.Range("L86").Value = 5224.5
$ws.Range("M86").Value = -424.25
$ws.Range("N86").Value = -7470.5
$ws.Range("H89").Value = 2466.5625
$ws.Range("I89").Value = 1547.25
$ws.Range("J89").Value = 5224.5
$ws.Range("K89").Value = 7736.25
$ws.Range("L89").Value = 26122.5
$ws.Range("M89").Value = -2120.25
$ws.Range("N89").Value = -37354.5
$ws.Range("H105").Value = 21197.666
$ws.Range("I105").Value = 51294.5
$ws.Range("J105").Value = 6149.25
$ws.Range("K105").Value = 51294.5
$ws.Range("L105").Value = 6149.25
$ws.Range("M105").Value = -49547.5
$ws.Range("N105").Value = -9643.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2959.8774
$ws.Range("I31").Value = 2592.1072
$ws.Range("J31").Value = 3450.238
$ws.Range("K31").Value = 2592.1072
$ws.Range("L31").Value = 3450.238
$ws.Range("M31").Value = -2297.1072
$ws.Range("N31").Value = -4040.238
$ws.Range("H34").Value = 2959.8774
$ws.Range("I34").Value = 2592.1072
$ws.Range("J34").Value = 3450.238
$ws.Range("K34").Value = 2592.1072
$ws.Range("L34").Value = 3450.238
$ws.Range("M34").Value = -2390.1072
$ws.Range("N34").Value = -3854.238
$ws.Range("H58").Value = 1401.305
$ws.Range("I58").Value = 1308.2683
$ws.Range("K58").Value = 1308.2683
$ws.Range("M58").Value = -1105.2683
$ws.Range("H109").Value = 64164.668
$ws.Range("H132").Value = 24190.69
$ws.Range("I132").Value = 28018.475
$ws.Range("K132").Value = 84055.42499999999
$ws.Range("M132").Value = -81525.42499999999
$ws.Range("H136").Value = 1401.305
$ws.Range("I136").Value = 1308.2683
$ws.Range("K136").Value = 3924.8049
$ws.Range("M136").Value = -1374.8049

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 235.09091
$ws.Range("I6").Value = 235.09091
$ws.Range("K6").Value = 705.27273
$ws.Range("M6").Value = -592.27273
$ws.Range("H8").Value = 734
$ws.Range("I8").Value = 734
$ws.Range("K8").Value = 2202
$ws.Range("M8").Value = -2063
$ws.Range("H94").Value = 8649.904
$ws.Range("J94").Value = 10205.177
$ws.Range("L94").Value = 30615.531
$ws.Range("N94").Value = -31967.531
$ws.Range("H98").Value = 358.3
$ws.Range("I98").Value = 188
$ws.Range("J98").Value = 377.22223
$ws.Range("K98").Value = 564
$ws.Range("L98").Value = 1131.66669
$ws.Range("M98").Value = 934
$ws.Range("N98").Value = -4127.66669
$ws.Range("H133").Value = 4499.3335
$ws.Range("I133").Value = 4499.3335
$ws.Range("K133").Value = 13498.0005
$ws.Range("M133").Value = -8438.000499999998
$ws.Range("H140").Value = 25644338
$ws.Range("I140").Value = 47621984
$ws.Range("J140").Value = 3749.8333
$ws.Range("K140").Value = 142865952
$ws.Range("L140").Value = 11249.4999
$ws.Range("M140").Value = -142860772
$ws.Range("N140").Value = -21609.4999
$ws.Range("H141").Value = 6104.5
$ws.Range("I141").Value = 6104.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 18313.5
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -13133.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1999.6666
$ws.Range("I113").Value = 1999.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1999.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 170.3334
$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -43494
$ws.Range("H122").Value = 2476.5945
$ws.Range("I122").Value = 2511.5667
$ws.Range("J122").Value = 2326.7144
$ws.Range("K122").Value = 7534.7001
$ws.Range("L122").Value = 6980.1432
$ws.Range("M122").Value = -5084.7001
$ws.Range("N122").Value = -11880.1432
$ws.Range("H132").Value = 992821.4399999999
$ws.Range("I132").Value = 7921.278
$ws.Range("K132").Value = 23763.834
$ws.Range("M132").Value = -21233.834

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2279.75
$ws.Range("I61").Value = 2102.037
$ws.Range("J61").Value = 3239.4
$ws.Range("K61").Value = 2102.037
$ws.Range("L61").Value = 3239.4
$ws.Range("M61").Value = -1900.037
$ws.Range("N61").Value = -3643.4
$ws.Range("H87").Value = 500012260
$ws.Range("J87").Value = 500012260
$ws.Range("L87").Value = 500012260
$ws.Range("N87").Value = -500014506
$ws.Range("H90").Value = 500012260
$ws.Range("J90").Value = 500012260
$ws.Range("L90").Value = 1500036780
$ws.Range("N90").Value = -1500048012
$ws.Range("H93").Value = 1436.1666
$ws.Range("I93").Value = 1436.1666
$ws.Range("K93").Value = 1436.1666
$ws.Range("M93").Value = -188.1666
$ws.Range("H113").Value = 2279.75
$ws.Range("I113").Value = 2102.037
$ws.Range("J113").Value = 3239.4
$ws.Range("K113").Value = 2102.037
$ws.Range("L113").Value = 3239.4
$ws.Range("M113").Value = 67.96300000000019
$ws.Range("N113").Value = -7579.4
$ws.Range("H119").Value = 45420
$ws.Range("J119").Value = 45420
$ws.Range("L119").Value = 45420
$ws.Range("N119").Value = -55096
$ws.Range("H122").Value = 3429.348
$ws.Range("I122").Value = 3039.8147
$ws.Range("J122").Value = 3982.8948
$ws.Range("K122").Value = 9119.444100000001
$ws.Range("L122").Value = 11948.6844
$ws.Range("M122").Value = -6669.444100000001
$ws.Range("N122").Value = -16848.6844

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 22249.5
$ws.Range("J4").Value = 20699
$ws.Range("L4").Value = 20699
$ws.Range("N4").Value = -20925
$ws.Range("H51").Value = 4199.4
$ws.Range("I51").Value = 1499.25
$ws.Range("K51").Value = 1499.25
$ws.Range("M51").Value = -989.25
$ws.Range("H119").Value = 38499.5
$ws.Range("J119").Value = 38499.5
$ws.Range("L119").Value = 38499.5
$ws.Range("N119").Value = -48175.5
$ws.Range("H136").Value = 21098.176
$ws.Range("I136").Value = 27587.352
$ws.Range("K136").Value = 82762.056
$ws.Range("M136").Value = -80212.056
